$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (losing formatting like trailing zeros).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '61.502.62'
$ws.Range('E2').Value = '  +0.98%  '
$ws.Range('D3').Value = '3.380.81'
$ws.Range('E3').Value = '  +0.60%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '575.11'
$ws.Range('E5').Value = '  +0.89%  '
$ws.Range('D6').Value = '136.85'
$ws.Range('E6').Value = '  +1.24%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.380.94'
$ws.Range('E8').Value = '  +0.60%  '
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('D10').Value = '7.50'
$ws.Range('E10').Value = '  -1.19%  '
$ws.Range('E11').Value = '  +2.35%  '
$ws.Range('D12').Value = '0.390'
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('D13').Value = '3.956.96'
$ws.Range('E13').Value = '  +0.55%  '
$ws.Range('E14').Value = '  +2.48%  '
$ws.Range('E15').Value = '  +2.94%  '
$ws.Range('D16').Value = '25.97'
$ws.Range('E16').Value = '  +3.85%  '
$ws.Range('D17').Value = '3.379.97'
$ws.Range('E17').Value = '  +0.51%  '
$ws.Range('D18').Value = '61.568.28'
$ws.Range('E18').Value = '  +0.96%  '
$ws.Range('D19').Value = '14.09'
$ws.Range('E19').Value = '  +1.05%  '
$ws.Range('D20').Value = '5.86'
$ws.Range('E20').Value = '  +1.39%  '
$ws.Range('D21').Value = '9.39'
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').Value = '376.69'
$ws.Range('E22').Value = '  +0.75%  '
$ws.Range('D23').Value = '0.557'
$ws.Range('E23').Value = '  -2.77%  '
$ws.Range('D24').Value = '3.520.51'
$ws.Range('E24').Value = '  +0.74%  '
$ws.Range('D25').Value = '0.999'
$ws.Range('E26').Value = '  +8.39%  '
$ws.Range('D27').Value = '71.36'
$ws.Range('E27').Value = '  +1.04%  '
$ws.Range('E28').Value = '  +5.30%  '
$ws.Range('D29').Value = '7.49'
$ws.Range('E29').Value = '  -2.60%  '
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '8.26'
$ws.Range('E31').Value = '  +2.38%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').Value = '0.161'
$ws.Range('E32').Value = '  +4.86%  '
$ws.Range('E33').Value = '  +2.05%  '
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').Value = '23.46'
$ws.Range('E35').Value = '  +0.20%  '
$ws.Range('D36').Value = '5.28'
$ws.Range('E36').Value = '  -4.55%  '
$ws.Range('D37').Value = '6.83'
$ws.Range('E37').Value = '  -0.89%  '
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('D39').Value = '165.30'
$ws.Range('E39').Value = '  +1.10%  '
$ws.Range('E40').Value = '  -1.54%  '
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('D42').Value = '0.775'
$ws.Range('E42').Value = '  +2.43%  '
$ws.Range('E43').Value = '  +7.81%  '
$ws.Range('E44').Value = '  +1.63%  '
$ws.Range('E45').Value = '  +0.81%  '
$ws.Range('D46').Value = '41.52'
$ws.Range('E46').Value = '  +0.15%  '
$ws.Range('D47').Value = '24.82'
$ws.Range('E47').Value = '  +8.71%  '
$ws.Range('D48').Value = '6.83'
$ws.Range('E48').Value = '  -1.51%  '
$ws.Range('D49').Value = '22.83'
$ws.Range('E49').Value = '  -0.79%  '
$ws.Range('D50').Value = '2.345.08'
$ws.Range('E50').Value = '  +4.60%  '
$ws.Range('D51').Value = '0.0262'
$ws.Range('E51').Value = '  +1.80%  '
